# Automated monthly refresh (2025-12-01):
#  - "VENTAS POR GRUPO": clears this months settled amounts that have now rolled off
#    the open/pending tracker, and resets the corresponding "<n> de 36" counters.
#  - "VENTA MENSUAL": rolls the 4-month rolling window forward by one month
#    (agosto -> septiembre, septiembre -> octubre, octubre -> noviembre,
#    noviembre -> diciembre), so each column's values shift right and a fresh,
#    still-empty "diciembre" column appears on the right.

$wb = $excel.ActiveWorkbook

# Sheet1 "VENTAS POR GRUPO": zero-out this period's closed/settled order amounts
# (by-product-category pivot) for the affected asesor/cliente rows, and reset the
# matching "<n> de 36" progress counters in the totals row to "0 de 36".
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M2").Value = 0

$ws1.Range("L3").Value = 0
$ws1.Range("M3").Value = 0

$ws1.Range("L5").Value = 0
$ws1.Range("M5").Value = 0

$ws1.Range("H9").Value = 0
$ws1.Range("I9").Value = 0
$ws1.Range("M9").Value = 0

$ws1.Range("M13").Value = 0

$ws1.Range("D16").Value = 0
$ws1.Range("M16").Value = 0

$ws1.Range("H21").Value = 0
$ws1.Range("L21").Value = 0
$ws1.Range("M21").Value = 0

$ws1.Range("D24").Value = 0

$ws1.Range("H30").Value = 0
$ws1.Range("I30").Value = 0

$ws1.Range("M33").Value = 0

$ws1.Range("M34").Value = 0

$ws1.Range("M35").Value = 0

$ws1.Range("D38").Value = "0 de 36"
$ws1.Range("H38").Value = "0 de 36"
$ws1.Range("I38").Value = "0 de 36"
$ws1.Range("L38").Value = "0 de 36"
$ws1.Range("M38").Value = "0 de 36"

# Sheet2 "VENTA MENSUAL": shift monthly columns right by one (agosto data rolls off,
# a new "diciembre" column is introduced with 0s), update widths and headers, and update
# shifted data values + the totals row.
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths (character units): C 14->16, D 16->14, E 14->15, F/G unchanged.
# The stored OOXML <col> width is ColumnWidth + 5/6, so compensate by subtracting 5/6.
$ws2.Columns.Item(3).ColumnWidth = 16 - (5/6)
$ws2.Columns.Item(4).ColumnWidth = 14 - (5/6)
$ws2.Columns.Item(5).ColumnWidth = 15 - (5/6)

# Header row: month names shift right by one month.
$ws2.Range("C1").Value = "septiembre"
$ws2.Range("D1").Value = "octubre"
$ws2.Range("E1").Value = "noviembre"
$ws2.Range("F1").Value = "diciembre"

$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 6045.88
$ws2.Range("E2").Value = 784.0599999999999
$ws2.Range("F2").Value = 0

$ws2.Range("C3").Value = 6231.33
$ws2.Range("D3").Value = 1326.66
$ws2.Range("E3").Value = 1496.52
$ws2.Range("F3").Value = 0

$ws2.Range("C4").Value = 687.03
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 0

$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 476.59
$ws2.Range("E5").Value = 286
$ws2.Range("F5").Value = 0

$ws2.Range("C6").Value = 0
$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 0

$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 0
$ws2.Range("E7").Value = 0
$ws2.Range("F7").Value = 0

$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 0
$ws2.Range("E8").Value = 0
$ws2.Range("F8").Value = 0

$ws2.Range("C9").Value = 3558.2
$ws2.Range("D9").Value = 2601.5
$ws2.Range("E9").Value = 3401.69
$ws2.Range("F9").Value = 0

$ws2.Range("C10").Value = 0
$ws2.Range("D10").Value = 0
$ws2.Range("E10").Value = 0
$ws2.Range("F10").Value = 0

$ws2.Range("C11").Value = 226.8
$ws2.Range("D11").Value = -309.07
$ws2.Range("E11").Value = 0
$ws2.Range("F11").Value = 0

$ws2.Range("C12").Value = 0
$ws2.Range("D12").Value = 0
$ws2.Range("E12").Value = 0
$ws2.Range("F12").Value = 0

$ws2.Range("C13").Value = 998.71
$ws2.Range("D13").Value = 1314.1
$ws2.Range("E13").Value = 66.68000000000001
$ws2.Range("F13").Value = 0

$ws2.Range("C14").Value = 0
$ws2.Range("D14").Value = 0
$ws2.Range("E14").Value = 0
$ws2.Range("F14").Value = 0

$ws2.Range("C15").Value = 0
$ws2.Range("D15").Value = 0
$ws2.Range("E15").Value = 0
$ws2.Range("F15").Value = 0

$ws2.Range("C16").Value = 10853.08
$ws2.Range("D16").Value = 350.18
$ws2.Range("E16").Value = -5874.77
$ws2.Range("F16").Value = 0

$ws2.Range("C17").Value = 0
$ws2.Range("D17").Value = 0
$ws2.Range("E17").Value = 0
$ws2.Range("F17").Value = 0

$ws2.Range("C18").Value = 0
$ws2.Range("D18").Value = 0
$ws2.Range("E18").Value = 0
$ws2.Range("F18").Value = 0

$ws2.Range("C19").Value = 0
$ws2.Range("D19").Value = 0
$ws2.Range("E19").Value = 0
$ws2.Range("F19").Value = 0

$ws2.Range("C20").Value = 4077.46
$ws2.Range("D20").Value = 0
$ws2.Range("E20").Value = 0
$ws2.Range("F20").Value = 0

$ws2.Range("C21").Value = 4354.56
$ws2.Range("D21").Value = 12049.42
$ws2.Range("E21").Value = 3779.22
$ws2.Range("F21").Value = 0

$ws2.Range("C22").Value = 0
$ws2.Range("D22").Value = 5015.36
$ws2.Range("E22").Value = 0
$ws2.Range("F22").Value = 0

$ws2.Range("C23").Value = 0
$ws2.Range("D23").Value = 0
$ws2.Range("E23").Value = 0
$ws2.Range("F23").Value = 0

$ws2.Range("C24").Value = 0
$ws2.Range("D24").Value = 5179.53
$ws2.Range("E24").Value = 366.34
$ws2.Range("F24").Value = 0

$ws2.Range("C25").Value = 0
$ws2.Range("D25").Value = 0
$ws2.Range("E25").Value = 0
$ws2.Range("F25").Value = 0

$ws2.Range("C26").Value = 0
$ws2.Range("D26").Value = 0
$ws2.Range("E26").Value = 0
$ws2.Range("F26").Value = 0

$ws2.Range("C27").Value = 0
$ws2.Range("D27").Value = 6777.81
$ws2.Range("E27").Value = 0
$ws2.Range("F27").Value = 0

$ws2.Range("C28").Value = 0
$ws2.Range("D28").Value = 0
$ws2.Range("E28").Value = 0
$ws2.Range("F28").Value = 0

$ws2.Range("C29").Value = 0
$ws2.Range("D29").Value = 0
$ws2.Range("E29").Value = 0
$ws2.Range("F29").Value = 0

$ws2.Range("C30").Value = 948.92
$ws2.Range("D30").Value = 259.58
$ws2.Range("E30").Value = 998.1
$ws2.Range("F30").Value = 0

$ws2.Range("C31").Value = 1831.68
$ws2.Range("D31").Value = 0
$ws2.Range("E31").Value = 0
$ws2.Range("F31").Value = 0

$ws2.Range("C32").Value = 0
$ws2.Range("D32").Value = 0
$ws2.Range("E32").Value = 0
$ws2.Range("F32").Value = 0

$ws2.Range("C33").Value = 0
$ws2.Range("D33").Value = 2536.39
$ws2.Range("E33").Value = 557.5599999999999
$ws2.Range("F33").Value = 0

$ws2.Range("C34").Value = 0
$ws2.Range("D34").Value = 0
$ws2.Range("E34").Value = 59.02
$ws2.Range("F34").Value = 0

$ws2.Range("C35").Value = 0
$ws2.Range("D35").Value = 0
$ws2.Range("E35").Value = 2350.86
$ws2.Range("F35").Value = 0

$ws2.Range("C36").Value = 0
$ws2.Range("D36").Value = 0
$ws2.Range("E36").Value = 0
$ws2.Range("F36").Value = 0

$ws2.Range("C37").Value = 5238.25
$ws2.Range("D37").Value = 1758.38
$ws2.Range("E37").Value = 0
$ws2.Range("F37").Value = 0

$ws2.Range("C38").Value = 39006.02
$ws2.Range("D38").Value = 45382.31
$ws2.Range("E38").Value = 8271.279999999999
$ws2.Range("F38").Value = 0
